$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173, shifting existing rows 173-319 down to 174-320.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new data point.
$ws.Cells.Item(173, 1).Value = 8
$ws.Cells.Item(173, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44827
$ws.Cells.Item(173, 5).Value = 4
$ws.Cells.Item(173, 6).Value = 100112012
$ws.Cells.Item(173, 7).Value = "Espinaca"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 2800
$ws.Cells.Item(173, 11).Value = 450
$ws.Cells.Item(173, 12).Value = 500
$ws.Cells.Item(173, 13).Value = 475
$ws.Cells.Item(173, 14).Value = '$/atado 300 a 500 gramos'
$ws.Cells.Item(173, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(173, 16).Value = 950
$ws.Cells.Item(173, 17).Value = 0.5
$ws.Cells.Item(173, 18).Value = "Hortaliza"
